# Commit: "changes in action class openURL"
#
# The automation run's "Results" column (column H, shared strings
# "Pass"/"Fail") is reset on three sheets:
#   - VerifyCSVForExistingVersion: rows 8,10,12,14,16 flip Pass -> Fail,
#     rows 17-44 have their result cleared entirely.
#   - VerifyCSVForNewVersion: rows 5-22 have their result cleared entirely.
#   - VerifyEventAPI: rows 5-62 have their result cleared entirely.

$wb = $excel.ActiveWorkbook

# --- VerifyCSVForExistingVersion -------------------------------------
$ws1 = $wb.Worksheets.Item("VerifyCSVForExistingVersion")

$failRows1 = @(8, 10, 12, 14, 16)
foreach ($r in $failRows1) {
    $ws1.Cells.Item($r, 8).Value = "Fail"
}

for ($r = 17; $r -le 44; $r++) {
    $ws1.Cells.Item($r, 8).ClearContents()
}

# --- VerifyCSVForNewVersion -------------------------------------------
$ws2 = $wb.Worksheets.Item("VerifyCSVForNewVersion")

for ($r = 5; $r -le 22; $r++) {
    $ws2.Cells.Item($r, 8).ClearContents()
}

# --- VerifyEventAPI -----------------------------------------------------
$ws3 = $wb.Worksheets.Item("VerifyEventAPI")

for ($r = 5; $r -le 62; $r++) {
    $ws3.Cells.Item($r, 8).ClearContents()
}
